$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10 (test #8): fill in date, time, condition and result -----------
# Set the cell values first; PasteSpecial(formats) afterwards restores the
# exact cell styles used by the other data rows (assigning .Value alone
# resets formatting to the default style).
$ws.Range("B10").Value = 43223
$ws.Range("C10").Value = 0.60277777777777775
$ws.Range("D10").Value = "- Without GPS module, time is completely synchronized from internet NTP Time server (more specific in result picture) and from each other`n- Lora transmitter is connected to laptop via usb uart CP2102 and transmit 26 packages each reset`n- Test indoor, room E6.1, 2 GPS with next by the window`n- 2 LoRa receiver is put side by side"
$ws.Range("E10").Value = "-Result in file indoor8"

# B10 needs the same date-formatted style used by the other "Date" cells.
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# D10 needs the same wrapped "Condition" style used by the other rows.
$ws.Range("D3").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# E10 ("Result") gets the same base style as the condition column, but
# without word-wrap (new cell style in the workbook).
$ws.Range("D3").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").WrapText = $false

# Row 10 grows to fit the new multi-line text.
$ws.Rows.Item(10).RowHeight = 159.75

# Move the active selection, matching the saved view state.
$ws.Range("E11").Select()
